# Auto-generated edit script for 'Fix heat rate modeling syntax'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76155.74719999998
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 9285.872463556781
$ws.Range("E2").Value = 2365
$ws.Range("F2").Value = 25216.55341308785

$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 104

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 161
$ws.Range("D4").Value = 0

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 20.8
$ws.Range("H2").Value = 41.6
$ws.Range("I2").Value = 52
$ws.Range("J2").Value = 57.60404040404041
$ws.Range("K2").Value = 72.8
$ws.Range("L2").Value = 83.2
$ws.Range("M2").Value = 93.59999999999999
$ws.Range("N2").Value = 104
$ws.Range("O2").Value = 93.59999999999999
$ws.Range("P2").Value = 83.2
$ws.Range("Q2").Value = 72.8
$ws.Range("R2").Value = 52
$ws.Range("S2").Value = 31.2
$ws.Range("T2").Value = 20.8

$ws.Range("I3").Value = 16.25311702887463
$ws.Range("J3").Value = 62.4
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 93.59999999999999
$ws.Range("M3").Value = 23.4
$ws.Range("N3").Value = 83.2
$ws.Range("O3").Value = 72.8
$ws.Range("Q3").Value = 52
$ws.Range("R3").Value = 31.2
$ws.Range("S3").Value = 20.8

$ws.Range("J4").Value = 10.4
$ws.Range("K4").Value = 41.6
$ws.Range("L4").Value = 72.8
$ws.Range("M4").Value = 23.4
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 72.8
$ws.Range("P4").Value = 37.38312417100297
$ws.Range("R4").Value = 10.4

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 13
$ws.Range("H2").Value = 28.6
$ws.Range("I2").Value = 20.8
$ws.Range("J2").Value = 18.60404040404041
$ws.Range("K2").Value = 46.8
$ws.Range("L2").Value = 62.4
$ws.Range("M2").Value = 70.2
$ws.Range("N2").Value = 78
$ws.Range("O2").Value = 62.4
$ws.Range("P2").Value = 54.6
$ws.Range("Q2").Value = 46.8
$ws.Range("R2").Value = 18.2

$ws.Range("I3").Value = 16.25311702887463
$ws.Range("J3").Value = 62.4
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 93.59999999999999
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 57.2
$ws.Range("O3").Value = 72.8
$ws.Range("Q3").Value = 26
$ws.Range("R3").Value = 31.2

$ws.Range("J4").Value = 10.4
$ws.Range("K4").Value = 41.6
$ws.Range("L4").Value = 72.8
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 72.8
$ws.Range("P4").Value = 37.38312417100297
$ws.Range("R4").Value = 10.4

$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("S2").Value = 10.4
$ws.Range("T2").Value = 25.14799999999994

$ws.Range("S3").Value = 20.8

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 187.8909090909091
$ws.Range("C2").Value = 168.1939393939394
$ws.Range("D2").Value = 155.0626262626263
$ws.Range("E2").Value = 141.9313131313131
$ws.Range("F2").Value = 128.8
$ws.Range("G2").Value = 141.67
$ws.Range("H2").Value = 169.984
$ws.Range("I2").Value = 190.576
$ws.Range("J2").Value = 208.994
$ws.Range("K2").Value = 255.326
$ws.Range("L2").Value = 317.102
$ws.Range("M2").Value = 386.6
$ws.Range("N2").Value = 463.82
$ws.Range("O2").Value = 525.596
$ws.Range("P2").Value = 579.65
$ws.Range("Q2").Value = 625.982
$ws.Range("R2").Value = 644
$ws.Range("S2").Value = 633.4949494949495
$ws.Range("T2").Value = 608.0929292929294
$ws.Range("U2").Value = 489.9111111111111
$ws.Range("V2").Value = 391.4262626262627
$ws.Range("W2").Value = 312.6383838383838
$ws.Range("X2").Value = 260.1131313131313
$ws.Range("Y2").Value = 220.7191919191919

$ws.Range("B3").Value = 181.3252525252525
$ws.Range("C3").Value = 161.6282828282828
$ws.Range("D3").Value = 148.4969696969697
$ws.Range("E3").Value = 148.4969696969697
$ws.Range("F3").Value = 148.4969696969697
$ws.Range("G3").Value = 128.8
$ws.Range("H3").Value = 128.8
$ws.Range("I3").Value = 144.8905858585859
$ws.Range("J3").Value = 206.6665858585859
$ws.Range("K3").Value = 206.6665858585859
$ws.Range("L3").Value = 299.3305858585859
$ws.Range("M3").Value = 299.3305858585859
$ws.Range("N3").Value = 355.9585858585859
$ws.Range("O3").Value = 428.0305858585859
$ws.Range("P3").Value = 428.0305858585859
$ws.Range("Q3").Value = 453.7705858585859
$ws.Range("R3").Value = 484.6585858585859
$ws.Range("S3").Value = 463.6484848484848
$ws.Range("T3").Value = 332.3353535353535
$ws.Range("U3").Value = 332.3353535353535
$ws.Range("V3").Value = 332.3353535353535
$ws.Range("W3").Value = 253.5474747474748
$ws.Range("X3").Value = 253.5474747474748
$ws.Range("Y3").Value = 214.1535353535353

$ws.Range("B4").Value = 168.1939393939394
$ws.Range("C4").Value = 148.4969696969697
$ws.Range("D4").Value = 148.4969696969697
$ws.Range("E4").Value = 148.4969696969697
$ws.Range("F4").Value = 148.4969696969697
$ws.Range("G4").Value = 128.8
$ws.Range("H4").Value = 128.8
$ws.Range("I4").Value = 128.8
$ws.Range("J4").Value = 139.096
$ws.Range("K4").Value = 180.28
$ws.Range("L4").Value = 252.352
$ws.Range("M4").Value = 252.352
$ws.Range("N4").Value = 252.352
$ws.Range("O4").Value = 324.424
$ws.Range("P4").Value = 361.4332929292929
$ws.Range("Q4").Value = 361.4332929292929
$ws.Range("R4").Value = 371.7292929292929
$ws.Range("S4").Value = 371.7292929292929
$ws.Range("T4").Value = 240.4161616161616
$ws.Range("U4").Value = 240.4161616161616
$ws.Range("V4").Value = 240.4161616161616
$ws.Range("W4").Value = 240.4161616161616
$ws.Range("X4").Value = 240.4161616161616
$ws.Range("Y4").Value = 201.0222222222222

$ws = $wb.Worksheets.Item("Feed in from Type 1")
$ws.Range("T2").Value = 6.052000000000064
